$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 23, shifting existing rows 23-110 down to rows 26-113
$ws.Range("A23:A25").EntireRow.Insert()

# New row 23 data
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C23").Value = 'Arica y Parinacota'
$ws.Range("D23").Value = 45177
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112045
$ws.Range("G23").Value = 'Zapallo'
$ws.Range("H23").Value = 'Camote'
$ws.Range("I23").Value = '1a nueva(o)'
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 650
$ws.Range("L23").Value = 670
$ws.Range("M23").Value = 663
$ws.Range("N23").Value = '$/kilo (volumen en unidades)'
$ws.Range("O23").Value = 'Perú'
$ws.Range("P23").Value = 663
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 'Hortaliza'

# New row 24 data
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C24").Value = 'Arica y Parinacota'
$ws.Range("D24").Value = 45177
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 100112045
$ws.Range("G24").Value = 'Zapallo'
$ws.Range("H24").Value = 'Camote'
$ws.Range("I24").Value = '2a nueva(o)'
$ws.Range("J24").Value = 630
$ws.Range("K24").Value = 630
$ws.Range("L24").Value = 640
$ws.Range("M24").Value = 636
$ws.Range("N24").Value = '$/kilo (volumen en unidades)'
$ws.Range("O24").Value = 'Perú'
$ws.Range("P24").Value = 636
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 'Hortaliza'

# New row 25 data
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C25").Value = 'Arica y Parinacota'
$ws.Range("D25").Value = 45177
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 100112045
$ws.Range("G25").Value = 'Zapallo'
$ws.Range("H25").Value = 'Camote'
$ws.Range("I25").Value = '3a nueva (o)'
$ws.Range("J25").Value = 330
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 620
$ws.Range("M25").Value = 611
$ws.Range("N25").Value = '$/kilo (volumen en unidades)'
$ws.Range("O25").Value = 'Perú'
$ws.Range("P25").Value = 611
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 'Hortaliza'
